# Generate Report for Handoff
# Update status "In Translation" -> "Ready for handoff" and bump the
# handoff timestamps on the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (B2, C2) and handoff date (D2)
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-03-22 20:38:50"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (E2)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-22 20:38:46"

# de-de sheet: Status (C2) and Latest Handoff Datetime (E2)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-22 20:38:50"
